# Update AgTests (F) and AgPosit (G) figures for rows 334-389
# Commit message: Updated: ut 30. 03. 2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F334").Value = 196712
$ws.Range("G334").Value = 3501
$ws.Range("F335").Value = 130699
$ws.Range("G335").Value = 2986
$ws.Range("F336").Value = 101774
$ws.Range("G336").Value = 3343
$ws.Range("F337").Value = 103934
$ws.Range("G337").Value = 2928
$ws.Range("F338").Value = 227110
$ws.Range("G338").Value = 3187
$ws.Range("F341").Value = 291721
$ws.Range("G341").Value = 3668
$ws.Range("F342").Value = 177429
$ws.Range("G342").Value = 3017
$ws.Range("F343").Value = 132789
$ws.Range("G343").Value = 2970
$ws.Range("F344").Value = 135235
$ws.Range("G344").Value = 2478
$ws.Range("F345").Value = 291964
$ws.Range("G345").Value = 3325
$ws.Range("F346").Value = 675100
$ws.Range("G346").Value = 4840
$ws.Range("F347").Value = 343972
$ws.Range("G347").Value = 2922
$ws.Range("F348").Value = 231849
$ws.Range("G348").Value = 3238
$ws.Range("F349").Value = 159100
$ws.Range("G349").Value = 2758
$ws.Range("F350").Value = 127336
$ws.Range("G350").Value = 2786
$ws.Range("F351").Value = 150645
$ws.Range("F352").Value = 307529
$ws.Range("G352").Value = 3548
$ws.Range("F353").Value = 724691
$ws.Range("G353").Value = 5292
$ws.Range("F355").Value = 222007
$ws.Range("G355").Value = 3448
$ws.Range("F356").Value = 160065
$ws.Range("G356").Value = 2877
$ws.Range("F357").Value = 138259
$ws.Range("G357").Value = 3025
$ws.Range("F358").Value = 157477
$ws.Range("G358").Value = 2601
$ws.Range("F359").Value = 321002
$ws.Range("G359").Value = 3348
$ws.Range("F360").Value = 748285
$ws.Range("G360").Value = 5130
$ws.Range("F362").Value = 228036
$ws.Range("G362").Value = 3173
$ws.Range("F363").Value = 187992
$ws.Range("G363").Value = 2763
$ws.Range("F364").Value = 167146
$ws.Range("G364").Value = 2460
$ws.Range("F365").Value = 183301
$ws.Range("G365").Value = 2390
$ws.Range("F366").Value = 338777
$ws.Range("F368").Value = 345789
$ws.Range("F369").Value = 233154
$ws.Range("G369").Value = 2589
$ws.Range("F370").Value = 181638
$ws.Range("G370").Value = 2030
$ws.Range("F371").Value = 158202
$ws.Range("G371").Value = 1947
$ws.Range("F372").Value = 177459
$ws.Range("F373").Value = 346542
$ws.Range("G373").Value = 2354
$ws.Range("F376").Value = 220444
$ws.Range("G376").Value = 2218
$ws.Range("F377").Value = 175651
$ws.Range("G377").Value = 1807
$ws.Range("F378").Value = 156298
$ws.Range("G378").Value = 1535
$ws.Range("F379").Value = 176657
$ws.Range("G379").Value = 1590
$ws.Range("F380").Value = 341653
$ws.Range("G380").Value = 1985
$ws.Range("F381").Value = 739429
$ws.Range("G381").Value = 2659
$ws.Range("F382").Value = 356173
$ws.Range("G382").Value = 1565
$ws.Range("F383").Value = 219448
$ws.Range("G383").Value = 1752
$ws.Range("F384").Value = 169290
$ws.Range("G384").Value = 1496
$ws.Range("F385").Value = 147224
$ws.Range("G385").Value = 1373
$ws.Range("F386").Value = 177219
$ws.Range("G386").Value = 1329
$ws.Range("F387").Value = 342237
$ws.Range("G387").Value = 1618
$ws.Range("F388").Value = 675358
$ws.Range("G388").Value = 2036
$ws.Range("F389").Value = 327888
$ws.Range("G389").Value = 1206
